# Refresh the cryptos list (prices + 1h volume %, plus a couple of rank swaps)
# as captured by the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.643.35"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.69%  "
$ws.Range("D3").Value = "'1.861.41"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.00%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'245.55"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.36%  "
$ws.Range("D6").Value = "'0.6989"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.07731"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.61%  "
$ws.Range("D9").Value = "'0.3069"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("D10").Value = "'23.68"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("D11").Value = "'0.07762"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").Value = "'5.165"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.58%  "
$ws.Range("D13").Value = "'1.856.64"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").Value = "'92.34"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.69%  "
$ws.Range("D15").Value = "'0.6930"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.35%  "
$ws.Range("D16").Value = "'6.557"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.36%  "
$ws.Range("D17").Value = "'29.623.47"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("D18").Value = "'0.000008363"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D19").Value = "'2.107.16"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.43%  "
$ws.Range("D20").Value = "'241.92"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("D22").Value = "'0.9999"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'7.613"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.19%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "'0.1505"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.44%  "
$ws.Range("D26").Value = "'8.916"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("D27").Value = "'159.60"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("D28").Value = "'18.31"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("D29").Value = "'1.532"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "'4.256"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("D31").Value = "'4.184"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("D32").Value = "'1.194"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("D33").Value = "'0.05103"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").Value = "'0.7833"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.70%  "
$ws.Range("D35").Value = "'1.902"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.33%  "
$ws.Range("D36").Value = "'1.157"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.80%  "
$ws.Range("D37").Value = "'2.687"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("D38").Value = "'1.324.83"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +10.83%  "
$ws.Range("D39").Value = "'0.01878"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("D40").Value = "'2.733"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.46%  "
$ws.Range("D41").Value = "'0.9608"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.02%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'106.45"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.857"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +13.34%  "
$ws.Range("D44").Value = "'1.0000"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'9.763"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.44%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000125"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.82%  "
$ws.Range("D47").Value = "'2.006.68"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("D48").Value = "'0.5215"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("D49").Value = "'1.788"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.13%  "
$ws.Range("D50").Value = "'64.49"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.69%  "
$ws.Range("D51").Value = "'7.006"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.20%  "
